# Generate Report for Handback
#
# 1. "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it
#    is used (Overview!B2:C2/B3:C3, zh-cn!C2/C3, de-de!C2/C3).
# 2. Fill in the "Latest Target File" / "Latest Handback File" columns (F/G)
#    for the zh-cn and de-de detail sheets, with hyperlinks matching the
#    source-file / handoff-file links already on each row.
# 3. Stamp the handback datetime into column H for both detail sheets
#    (zh-cn rows get one timestamp, de-de rows get another).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $newStatus
$ovw.Range("C2").Value = $newStatus
$ovw.Range("B3").Value = $newStatus
$ovw.Range("C3").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

# Row 2: target / handback file links mirror the source (A2) / handoff (D2) links
$zh.Range("F2").Value = "196f4342-cf7e-4c67-b105-f0f976a4b81c.md"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ac241f733090812ea89fc453bb837485783df82/e2e/196f4342-cf7e-4c67-b105-f0f976a4b81c.md", "", "", "196f4342-cf7e-4c67-b105-f0f976a4b81c.md") | Out-Null

$zh.Range("G2").Value = "196f4342-cf7e-4c67-b105-f0f976a4b81c.cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/400fe82060534da9234ba71bff3bf0203a6cad37/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/196f4342-cf7e-4c67-b105-f0f976a4b81c.cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7.zh-cn.xlf", "", "", "196f4342-cf7e-4c67-b105-f0f976a4b81c.cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7.zh-cn.xlf") | Out-Null

# Row 3
$zh.Range("F3").Value = "650b5967-7ade-4182-9ac8-804e3ebc3ae8.md"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ac241f733090812ea89fc453bb837485783df82/e2e/650b5967-7ade-4182-9ac8-804e3ebc3ae8.md", "", "", "650b5967-7ade-4182-9ac8-804e3ebc3ae8.md") | Out-Null

$zh.Range("G3").Value = "650b5967-7ade-4182-9ac8-804e3ebc3ae8.918378bc386f3928fb7cab8fb055ddcf83834c05.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/400fe82060534da9234ba71bff3bf0203a6cad37/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/650b5967-7ade-4182-9ac8-804e3ebc3ae8.918378bc386f3928fb7cab8fb055ddcf83834c05.zh-cn.xlf", "", "", "650b5967-7ade-4182-9ac8-804e3ebc3ae8.918378bc386f3928fb7cab8fb055ddcf83834c05.zh-cn.xlf") | Out-Null

# Handback datetime stamp (shared by both rows on this sheet)
$zh.Range("H2").Value = "2016-03-20 22:54:37"
$zh.Range("H3").Value = "2016-03-20 22:54:37"

# --- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# Row 2
$de.Range("F2").Value = "196f4342-cf7e-4c67-b105-f0f976a4b81c.md"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ac241f733090812ea89fc453bb837485783df82/e2e/196f4342-cf7e-4c67-b105-f0f976a4b81c.md", "", "", "196f4342-cf7e-4c67-b105-f0f976a4b81c.md") | Out-Null

$de.Range("G2").Value = "196f4342-cf7e-4c67-b105-f0f976a4b81c.cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/464f0f41e0a2d6eab20d971cf32e31766ea0a7d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/196f4342-cf7e-4c67-b105-f0f976a4b81c.cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7.de-de.xlf", "", "", "196f4342-cf7e-4c67-b105-f0f976a4b81c.cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7.de-de.xlf") | Out-Null

# Row 3
$de.Range("F3").Value = "650b5967-7ade-4182-9ac8-804e3ebc3ae8.md"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ac241f733090812ea89fc453bb837485783df82/e2e/650b5967-7ade-4182-9ac8-804e3ebc3ae8.md", "", "", "650b5967-7ade-4182-9ac8-804e3ebc3ae8.md") | Out-Null

$de.Range("G3").Value = "650b5967-7ade-4182-9ac8-804e3ebc3ae8.918378bc386f3928fb7cab8fb055ddcf83834c05.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/464f0f41e0a2d6eab20d971cf32e31766ea0a7d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/650b5967-7ade-4182-9ac8-804e3ebc3ae8.918378bc386f3928fb7cab8fb055ddcf83834c05.de-de.xlf", "", "", "650b5967-7ade-4182-9ac8-804e3ebc3ae8.918378bc386f3928fb7cab8fb055ddcf83834c05.de-de.xlf") | Out-Null

# Handback datetime stamp (shared by both rows on this sheet, distinct from zh-cn's)
$de.Range("H2").Value = "2016-03-20 22:54:42"
$de.Range("H3").Value = "2016-03-20 22:54:42"
